$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as Text, avoiding Excel's automatic
# number inference (and the float round-trip noise that introduces),
# then reset the style so no residual formatting is left on the cell.
function Set-TextCell {
    param($Sheet, [string]$Addr, [string]$Val)
    $r = $Sheet.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "37.147.46"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.079.11"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextCell $ws "D5" "253.67"
$ws.Range("E5").Value = "  +1.09%  "
Set-TextCell $ws "D6" "0.677"
$ws.Range("E6").Value = "  +1.81%  "
Set-TextCell $ws "D7" "59.25"
$ws.Range("E7").Value = "  +9.07%  "
Set-TextCell $ws "D9" "0.393"
$ws.Range("E9").Value = "  +4.59%  "
Set-TextCell $ws "D10" "61.54"
$ws.Range("E10").Value = "  -0.64%  "
Set-TextCell $ws "D11" "0.0803"
$ws.Range("E11").Value = "  +7.53%  "
Set-TextCell $ws "D12" "0.109"
$ws.Range("E12").Value = "  +2.48%  "
Set-TextCell $ws "D13" "16.39"
$ws.Range("E13").Value = "  +7.02%  "
$ws.Range("D14").Value = "2.381.99"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("E15").Value = "  -2.09%  "
Set-TextCell $ws "D16" "5.56"
$ws.Range("E16").Value = "  +7.32%  "
$ws.Range("D17").Value = "2.077.34"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "37.560.87"
$ws.Range("E18").Value = "  +0.51%  "
Set-TextCell $ws "D19" "15.79"
$ws.Range("E19").Value = "  +6.79%  "
Set-TextCell $ws "D20" "74.79"
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("D21").Value = "0.0₃0928"
$ws.Range("E21").Value = "  +8.95%  "
$ws.Range("E22").Value = "  +4.91%  "
Set-TextCell $ws "D23" "239.68"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -2.49%  "
Set-TextCell $ws "D26" "2.28"
$ws.Range("E26").Value = "  +13.93%  "
Set-TextCell $ws "D27" "169.90"
$ws.Range("E27").Value = "  -1.50%  "
Set-TextCell $ws "D28" "9.35"
$ws.Range("E28").Value = "  +0.74%  "
Set-TextCell $ws "D29" "20.42"
$ws.Range("E29").Value = "  -1.64%  "
Set-TextCell $ws "D30" "0.127"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("E31").Value = "  +6.41%  "
$ws.Range("E32").Value = "  +5.87%  "
Set-TextCell $ws "D33" "0.0636"
$ws.Range("E33").Value = "  +2.99%  "
Set-TextCell $ws "D34" "4.51"
$ws.Range("E34").Value = "  +8.75%  "
Set-TextCell $ws "D35" "0.0912"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("E36").Value = "  +0.05%  "
Set-TextCell $ws "D37" "2.32"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("B38").Value = "Cronos"
$ws.Range("C38").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws "D38" "0.116"
$ws.Range("E38").Value = "  +25.68%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D39" "1.77"
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D41" "0.0227"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws "D42" "17.86"
$ws.Range("E42").Value = "  -3.71%  "
Set-TextCell $ws "D43" "1.18"
$ws.Range("E43").Value = "  +0.69%  "
Set-TextCell $ws "D44" "99.41"
$ws.Range("E44").Value = "  +0.23%  "
Set-TextCell $ws "D45" "4.32"
$ws.Range("E45").Value = "  +7.33%  "
Set-TextCell $ws "D46" "2.84"
$ws.Range("E46").Value = "  +1.27%  "
Set-TextCell $ws "D47" "4.61"
$ws.Range("E47").Value = "  +14.73%  "
$ws.Range("E48").Value = "  +8.07%  "
$ws.Range("D49").Value = "1.309.81"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("E50").Value = "  -0.17%  "
Set-TextCell $ws "D51" "6.95"
$ws.Range("E51").Value = "  -0.95%  "
